$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (Volume/Number + report week dates)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/1/2024  Through  4/7/2024"

# ---------------------------------------------------------------------------
# Number formats reused from the existing style table so that writing values
# keeps the same cellXf/style ids as already used throughout the sheet.
# ---------------------------------------------------------------------------
$fmtCount = "#,##0"
$fmtPct   = '#,##0.0;"-"#,##0.0'

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = 2
$ws.Range("C15").NumberFormat = $fmtCount
$ws.Range("F15").Value = 2
$ws.Range("F15").NumberFormat = $fmtCount
$ws.Range("G15").Value2 = 3
$ws.Range("H15").Value2 = -33.333333333333
$ws.Range("I15").Value2 = 4
$ws.Range("K15").Value2 = -33.333333333333
$ws.Range("L15").Value2 = -66.666666666666
$ws.Range("M15").Value2 = -20
$ws.Range("N15").Value2 = -88.235294117647

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$ws.Range("D16").Value2 = 4
$ws.Range("E16").Value2 = -25
$ws.Range("G16").Value2 = 12
$ws.Range("H16").Value2 = -16.666666666666
$ws.Range("I16").Value2 = 45
$ws.Range("J16").Value2 = 44
$ws.Range("K16").Value2 = 2.272727272727
$ws.Range("L16").Value2 = -33.823529411764
$ws.Range("N16").Value2 = -92.227979274611

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Range("C17").Value2 = 3
$ws.Range("E17").Value2 = -40
$ws.Range("F17").Value2 = 11
$ws.Range("G17").Value2 = 27
$ws.Range("H17").Value2 = -59.259259259259
$ws.Range("I17").Value2 = 65
$ws.Range("J17").Value2 = 99
$ws.Range("K17").Value2 = -34.343434343434
$ws.Range("L17").Value2 = -34.343434343434
$ws.Range("M17").Value2 = -34.343434343434
$ws.Range("N17").Value2 = -68.899521531100

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
$ws.Range("C18").Value2 = 2
$ws.Range("E18").Value2 = 100
$ws.Range("F18").Value2 = 13
$ws.Range("G18").Value2 = 8
$ws.Range("H18").Value2 = 62.5
$ws.Range("I18").Value2 = 36
$ws.Range("J18").Value2 = 43
$ws.Range("K18").Value2 = -16.279069767441
$ws.Range("L18").Value2 = -21.739130434782
$ws.Range("M18").Value2 = -40.983606557377
$ws.Range("N18").Value2 = -95.577395577395

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Range("C19").Value2 = 10
$ws.Range("D19").Value2 = 6
$ws.Range("E19").Value2 = 66.666666666666
$ws.Range("F19").Value2 = 43
$ws.Range("G19").Value2 = 36
$ws.Range("H19").Value2 = 19.444444444444
$ws.Range("I19").Value2 = 124
$ws.Range("J19").Value2 = 154
$ws.Range("K19").Value2 = -19.480519480519
$ws.Range("L19").Value2 = -15.068493150684
$ws.Range("M19").Value2 = -16.778523489932
$ws.Range("N19").Value2 = -54.578754578754

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$ws.Range("C20").Value2 = 4
$ws.Range("D20").Value2 = 2
$ws.Range("E20").Value2 = 100
$ws.Range("F20").Value2 = 11
$ws.Range("G20").Value2 = 7
$ws.Range("H20").Value2 = 57.142857142857
$ws.Range("I20").Value2 = 33
$ws.Range("J20").Value2 = 26
$ws.Range("K20").Value2 = 26.923076923076
$ws.Range("L20").Value2 = -21.428571428571
$ws.Range("M20").Value2 = -42.105263157894
$ws.Range("N20").Value2 = -94.481605351170

# ---------------------------------------------------------------------------
# Row 21 (TOTAL row, bold styles - values only)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value2 = 24
$ws.Range("D21").Value2 = 18
$ws.Range("E21").Value2 = 33.333333333333
$ws.Range("F21").Value2 = 90
$ws.Range("G21").Value2 = 93
$ws.Range("H21").Value2 = -3.225806451612
$ws.Range("I21").Value2 = 307
$ws.Range("J21").Value2 = 373
$ws.Range("K21").Value2 = -17.694369973190
$ws.Range("L21").Value2 = -25.665859564164
$ws.Range("M21").Value2 = -35.908141962421
$ws.Range("N21").Value2 = -87.807783955520

# ---------------------------------------------------------------------------
# Row 22
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("F22").Value2 = 2
$ws.Range("H22").Value2 = 0

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("C24").Value2 = 26
$ws.Range("D24").Value2 = 38
$ws.Range("E24").Value2 = -31.578947368421
$ws.Range("F24").Value2 = 112
$ws.Range("H24").Value2 = -29.113924050632
$ws.Range("I24").Value2 = 513
$ws.Range("J24").Value2 = 464
$ws.Range("K24").Value2 = 10.560344827586
$ws.Range("L24").Value2 = 50.882352941176
$ws.Range("M24").Value2 = 65.483870967741

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("C25").Value2 = 20
$ws.Range("D25").Value2 = 20
$ws.Range("E25").Value2 = 0
$ws.Range("F25").Value2 = 78
$ws.Range("G25").Value2 = 96
$ws.Range("H25").Value2 = -18.75
$ws.Range("I25").Value2 = 312
$ws.Range("J25").Value2 = 267
$ws.Range("K25").Value2 = 16.853932584269
$ws.Range("L25").Value2 = 95

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$ws.Range("C26").Value2 = 6
$ws.Range("D26").Value2 = 19
$ws.Range("E26").Value2 = -68.421052631578
$ws.Range("F26").Value2 = 53
$ws.Range("G26").Value2 = 61
$ws.Range("H26").Value2 = -13.114754098360
$ws.Range("I26").Value2 = 172
$ws.Range("J26").Value2 = 157
$ws.Range("K26").Value2 = 9.554140127388
$ws.Range("L26").Value2 = 29.323308270676
$ws.Range("M26").Value2 = -10.880829015544

# ---------------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------------
$ws.Range("C27").Value2 = 3
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = $fmtCount
$ws.Range("E27").Value = 200
$ws.Range("E27").NumberFormat = $fmtPct
$ws.Range("F27").Value2 = 4
$ws.Range("H27").Value2 = -33.333333333333
$ws.Range("I27").Value2 = 10
$ws.Range("J27").Value2 = 13
$ws.Range("K27").Value2 = -23.076923076923
$ws.Range("L27").Value2 = -44.444444444444

# ---------------------------------------------------------------------------
# Row 28
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = $fmtCount
$ws.Range("D28").Value2 = 3
$ws.Range("E28").Value2 = -66.666666666666
$ws.Range("G28").Value2 = 9
$ws.Range("H28").Value2 = -33.333333333333
$ws.Range("I28").Value2 = 20
$ws.Range("J28").Value2 = 15
$ws.Range("K28").Value2 = 33.333333333333
$ws.Range("L28").Value2 = 33.333333333333

# ---------------------------------------------------------------------------
# Row 31
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))
$ws.Range("G31").Value2 = 3
$ws.Range("H31").Value2 = -66.666666666666
$ws.Range("I31").Value2 = 3
$ws.Range("K31").Value2 = -40
$ws.Range("L31").Value2 = 200

# ---------------------------------------------------------------------------
# Row 33
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("F33"))

# ---------------------------------------------------------------------------
# Column E width: now matches the narrower columns (C/D/F/G)
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
